# edit.ps1 - apply 'Documentacion de Diseno: descripcion capa fisica' commit
# Inserts the 'Capa Fisica' section body: ten new paragraphs describing
# ArchivoBloques / SerialBuffer, right after the existing 'Capa Física'
# Heading2 paragraph and before the 'Árbol B+' Heading2 paragraph.

$d = $word.ActiveDocument

function Get-ParaByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# Step 1: bulk insert the ten new paragraphs (plus a trailing 11th *empty*
# placeholder paragraph) right after the 'Capa Física' heading, in one single
# Find/Replace edit. The trailing empty paragraph is a scratch paragraph used
# in step 3 below to relocate the document's _GoBack bookmark to the true end
# of the inserted block (matching the original author's final cursor position).
$found = $d.Content.Find.Execute("Capa Física", $true, $false, $false, $false, $false, $true, 1, $false, "Capa Física^pLa capa física de este trabajo se centra alrededor de dos clases: ArchivoBloques y SerialBuffer.^pLa primera es la encargada de la interacción directa con el disco. La misma define interfaces para poder leer y escribir un archivo en función de bloques de un tamaño fijo y parametrizable.^pEl problema que surgió aquí es que independientemente de utilizar el modo normal de apertura para escribir de C++ (ios::out), el archivo era truncado y se perdía su contenido, quedando solamente el último bloque escrito. La solución para esto fue que si se quería modificar un bloque ya existente, se utilizara un archivo de trabajo temporal para pasar el contenido original, intercalar el bloque a modificar y luego copiar el resto del archivo. Posteriormente se elimina el archivo original y se renombra el de trabajo. Debido a que esto es muy costoso y generalmente muchas de las operaciones de escritura consisten en agregar un nuevo bloque al final, existe un método que abre el archivo en modo append (ios::app) y agrega el nuevo bloque. La lectura en cambio no presento problema alguno.^pSerialBuffer es la clase encargada de brindar los medios para poder persistir de manera ordenada los registros de cada estructura. La misma presenta dos métodos principales, pack y unpack.^pEl primero se encarga de agregar registros a un buffer de caracteres (cabe aclarar que se escogió por archivos de tipo binario para la persistencia de datos). Para poder recuperar la información a posteriori, antes de empaquetar cada registro carga en el buffer un prefijo de longitud para poder saber cuanto va a tener que recuperar del buffer a un objeto.^pEl segundo, hace el trabajo inverso y restaura la información de un buffer, que previamente se leyó desde el disco, a un objeto. Utiliza el prefijo de longitu para saber la cantidad de caracteres a pasar desde el buffer al objeto de destino.^pComo puede notarse, el buffer es la estructura fundamental y puede verse al mismo como una sucesión de registros así:^p[prefijoLongitud(unsigned short int), reg(registro genérico de longitud variable)]^pDonde el tamaño total queda comprendido dentro del tamaño del buffer.^pUna carácteristica de esta implementación es la relación uno a uno que debe haber entre los tamaños de bloque y de los buffer, ya que si los mismos no coincidieran generarían errores de segmentación a la hora de tratar datos en memoria.^p", 2)
if (-not $found) { throw "Could not find Capa Física heading paragraph" }

# Step 2: fix up paragraph formatting for each of the ten new paragraphs -
# strip the inherited Heading2 style/outline numbering and apply the plain
# body-text formatting used throughout the rest of the document.
$prefixes = @(
    "La capa física de este tr",
    "La primera es la encargad",
    "El problema que surgió aq",
    "SerialBuffer es la clase ",
    "El primero se encarga de ",
    "El segundo, hace el traba",
    "Como puede notarse, el bu",
    "[prefijoLongitud(unsigned",
    "Donde el tamaño total que",
    "Una carácteristica de est"
)

foreach ($prefix in $prefixes) {
    $p = Get-ParaByPrefix $prefix
    if ($p -eq $null) { throw "Could not locate inserted paragraph: $prefix" }
    $p.Range.ListFormat.RemoveNumbers()
    $p.Range.Style = "Normal"
    $p.Range.ParagraphFormat.Alignment = 3
}

# The eighth new paragraph (the '[prefijoLongitud...]' code sample) is
# centred and set in Arial Narrow instead of the body-text justification.
$codePara = Get-ParaByPrefix "[prefijoLongitud(unsigned"
if ($codePara -eq $null) { throw "Could not locate code sample paragraph" }
$codePara.Range.ParagraphFormat.Alignment = 1
$codePara.Range.Font.Name = "Arial Narrow"

# Step 3: relocate the _GoBack bookmark from the start of the trailing empty
# 11th paragraph to the true end of the tenth (last) new paragraph, by
# deleting the paragraph mark that separates them - this merges the now-empty
# 11th paragraph away and leaves the bookmark sitting right after the last run
# of paragraph ten, exactly where the original author's edit ended.
$lastPara = Get-ParaByPrefix "Una carácteristica de est"
if ($lastPara -eq $null) { throw "Could not locate last new paragraph" }
$emptyPara = $lastPara.Next()
$er = $emptyPara.Range.Duplicate
$markRange = $d.Range($er.Start - 1, $er.Start)
if ($markRange.Text -ne [string][char]13) { throw "Unexpected trailing content; aborting to avoid corrupting the document" }
$markRange.Delete()

Write-Output "done"
